$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price / 1h-volume refresh (scheduled GitHub Actions run).
# D-column values that parse as plain numbers get an apostrophe
# prefix so Excel keeps storing/displaying them as literal text
# (matching the sheet's existing text-formatted price strings)
# instead of silently converting them to floating point numbers.

$ws.Range("D2").Value = "63.169.15"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.557.10"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'568.24"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'146.63"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "2.553.29"
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "'27.52"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "3.008.22"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "63.100.08"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "2.553.60"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("D19").Value = "'11.43"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "'335.49"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "'6.80"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'65.26"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'1.64"
$ws.Range("E25").Value = "  +9.80%  "
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").Value = "'8.48"
$ws.Range("E27").Value = "  +4.81%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'1.47"
$ws.Range("E29").Value = "  +5.83%  "
$ws.Range("D30").Value = "'7.34"
$ws.Range("E30").Value = "  +7.45%  "
$ws.Range("D31").Value = "0.0₃0819"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("D33").Value = "'175.91"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("D35").Value = "'407.70"
$ws.Range("E35").Value = "  +9.03%  "
$ws.Range("D36").Value = "'0.398"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'18.99"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'39.32"
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("D43").Value = "'153.03"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").Value = "'21.05"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("D46").Value = "'0.607"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").Value = "'0.0527"
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("D48").Value = "'0.0961"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  +5.55%  "
$ws.Range("D50").Value = "'18.39"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("E51").Value = "  +0.81%  "
